# Update TPM-derived NATMI metrics (Gdf11-Acvr2b) to reflect the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3288063333333333
$ws.Range("H2").Value = 0.9864189999999999
$ws.Range("I2").Value = 0.05575527297994041
$ws.Range("J2").Value = 0.05575527297994041
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.315861666666667
$ws.Range("N2").Value = 3.947585
$ws.Range("O2").Value = 0.2754050739440597
$ws.Range("P2").Value = 0.2754050739440597
$ws.Range("Q2").Value = 0.4326636497905555
$ws.Range("R2").Value = 3.893972848115
$ws.Range("S2").Value = 0.01535528507781172
$ws.Range("T2").Value = 0.01535528507781172

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3288063333333333
$ws.Range("H3").Value = 0.9864189999999999
$ws.Range("I3").Value = 0.05575527297994041
$ws.Range("J3").Value = 0.05575527297994041
$ws.Range("O3").Value = 0.3040809095127364
$ws.Range("P3").Value = 0.3040809095127364
$ws.Range("Q3").Value = 0.4777136247247777
$ws.Range("R3").Value = 4.299422622522999
$ws.Range("S3").Value = 0.01695411411787118
$ws.Range("T3").Value = 0.01695411411787118

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3288063333333333
$ws.Range("H4").Value = 0.9864189999999999
$ws.Range("I4").Value = 0.05575527297994041
$ws.Range("J4").Value = 0.05575527297994041
$ws.Range("M4").Value = 2.009179666666667
$ws.Range("N4").Value = 6.027539
$ws.Range("O4").Value = 0.4205140165432039
$ws.Range("P4").Value = 0.4205140165432039
$ws.Range("Q4").Value = 0.6606309992045556
$ws.Range("R4").Value = 5.945678992841
$ws.Range("S4").Value = 0.02344587378425751
$ws.Range("T4").Value = 0.02344587378425751

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.3115445049245869
$ws.Range("J5").Value = 0.3115445049245869
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.315861666666667
$ws.Range("N5").Value = 3.947585
$ws.Range("O5").Value = 0.2754050739440597
$ws.Range("P5").Value = 0.2754050739440597
$ws.Range("Q5").Value = 2.417600620866111
$ws.Range("R5").Value = 21.758405587795
$ws.Range("S5").Value = 0.08580093741562132
$ws.Range("T5").Value = 0.08580093741562132

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.3115445049245869
$ws.Range("J6").Value = 0.3115445049245869
$ws.Range("O6").Value = 0.3040809095127364
$ws.Range("P6").Value = 0.3040809095127364
$ws.Range("S6").Value = 0.09473473641116356
$ws.Range("T6").Value = 0.09473473641116356

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.3115445049245869
$ws.Range("J7").Value = 0.3115445049245869
$ws.Range("M7").Value = 2.009179666666667
$ws.Range("N7").Value = 6.027539
$ws.Range("O7").Value = 0.4205140165432039
$ws.Range("P7").Value = 0.4205140165432039
$ws.Range("Q7").Value = 3.691416911528111
$ws.Range("R7").Value = 33.222752203753
$ws.Range("S7").Value = 0.131008831097802
$ws.Range("T7").Value = 0.131008831097802

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 3.731231666666667
$ws.Range("H8").Value = 11.193695
$ws.Range("I8").Value = 0.6327002220954728
$ws.Range("J8").Value = 0.6327002220954728
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.315861666666667
$ws.Range("N8").Value = 3.947585
$ws.Range("O8").Value = 0.2754050739440597
$ws.Range("P8").Value = 0.2754050739440597
$ws.Range("Q8").Value = 4.909784719619446
$ws.Range("R8").Value = 44.18806247657501
$ws.Range("S8").Value = 0.1742488514506267
$ws.Range("T8").Value = 0.1742488514506267

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 3.731231666666667
$ws.Range("H9").Value = 11.193695
$ws.Range("I9").Value = 0.6327002220954728
$ws.Range("J9").Value = 0.6327002220954728
$ws.Range("O9").Value = 0.3040809095127364
$ws.Range("P9").Value = 0.3040809095127364
$ws.Range("Q9").Value = 5.421003257757223
$ws.Range("R9").Value = 48.78902931981501
$ws.Range("S9").Value = 0.1923920589837017
$ws.Range("T9").Value = 0.1923920589837017

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 3.731231666666667
$ws.Range("H10").Value = 11.193695
$ws.Range("I10").Value = 0.6327002220954728
$ws.Range("J10").Value = 0.6327002220954728
$ws.Range("M10").Value = 2.009179666666667
$ws.Range("N10").Value = 6.027539
$ws.Range("O10").Value = 0.4205140165432039
$ws.Range("P10").Value = 0.4205140165432039
$ws.Range("Q10").Value = 7.496714796289447
$ws.Range("R10").Value = 67.470433166605
$ws.Range("S10").Value = 0.2660593116611444
$ws.Range("T10").Value = 0.2660593116611444

